$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.278.69'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '2.466.43'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.532'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('D9').Value = '2.463.09'
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').Value = '  +4.32%  '
$ws.Range('E11').Value = '  +2.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.21'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = '2.901.76'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').Value = '62.114.15'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '2.458.15'
$ws.Range('E18').Value = '  +1.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.14%  '
$ws.Range('E20').Value = '  +3.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '327.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '65.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '589.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.39%  '
$ws.Range('D29').Value = '2.575.89'
$ws.Range('E29').Value = '  +1.27%  '
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').Value = '0.0₃0948'
$ws.Range('E31').Value = '  +0.57%  '
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('E33').Value = '  -2.71%  '
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('E35').Value = '  -2.54%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.81'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.73%  '
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.374'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '150.98'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.35'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.32'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.41'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('D47').Value = '0.0₆0291'
$ws.Range('E47').Value = '  +23.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '143.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.604'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.21%  '
